# "unify the conception of DataNode, DataTable, Entity."
#
# - Rename the "Property*" sheets to "DataNode_*"
# - Rename the "Record_*" sheets (except Record_Building) to "DataTable_*"
# - Drop the now-obsolete "Record_Building" sheet entirely
# - Leave "Component" as-is

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$wb.Worksheets.Item("Property1").Name = "DataNode_1"
$wb.Worksheets.Item("Property2").Name = "DataNode_2"
$wb.Worksheets.Item("Record_Hero").Name = "DataTable_Hero"
$wb.Worksheets.Item("Record_Bag").Name = "DataTable_Bag"
$wb.Worksheets.Item("Record_CommPropertyValue").Name = "DataTable_CommPropertyValue"
$wb.Worksheets.Item("Record_Task").Name = "DataTable_Task"

$wb.Worksheets.Item("Record_Building").Delete() | Out-Null

$wb.Worksheets.Item("DataTable_Hero").Activate() | Out-Null
